# Sablona_za_nenaplacena_potrazivanja_primjer-1.xlsx
# - show invoice amounts on "Racuni" with 2 decimals + "HRK" suffix (comma as decimal point)
# - bump a couple of sample amounts so the new formatting is visible
# - move the active tab/selection around
# - fix the "Kupci" row numbering
# - flag G4 (text-stored number) as an ignored error

$wb = $excel.ActiveWorkbook

$wsZ = $wb.Worksheets.Item("Zaglavlje")
$wsK = $wb.Worksheets.Item("Kupci")
$wsR = $wb.Worksheets.Item("Racuni")

# ---------------------------------------------------------------------------
# Racuni: apply a new currency number format ("#,##0.00 HRK") to the
# Iznos/Pdv/Placeni iznos columns (D:F), header + data cells alike.
# ---------------------------------------------------------------------------
$hrkFormat = "#,##0.00"" ""HRK"

$wsR.Range("D1:F1").NumberFormat = $hrkFormat
$wsR.Range("D2:F4").NumberFormat = $hrkFormat

# Bump a couple of sample values so the 2-decimal formatting is visible
$wsR.Range("D2").Value = 10000.5
$wsR.Range("D3").Value = 20000.66
$wsR.Range("E3").Value = 200.66
$wsR.Range("F3").Value = 2000.66

# Column widths for D/E/F (13 / 12 / 19 chars) - E is a brand-new custom column
$wsR.Columns.Item(4).ColumnWidth = 12.1666666
$wsR.Columns.Item(5).ColumnWidth = 11.1666666
$wsR.Columns.Item(6).ColumnWidth = 18.1666666

# G4 holds a large OIB stored as text -- mark it as an intentionally ignored error
$wsR.Range("G4").Errors.Item(1).Ignore = $true

# ---------------------------------------------------------------------------
# Kupci: row 3 "redni broj" should read 1 (matches row 2) instead of 2
# ---------------------------------------------------------------------------
$wsK.Range("A3").Value = 1

# ---------------------------------------------------------------------------
# Selections per sheet (also determines which tab ends up "active")
# ---------------------------------------------------------------------------
$wsZ.Range("I2").Select() | Out-Null
$wsK.Range("D11").Select() | Out-Null
$wsR.Range("G6").Select() | Out-Null

# Racuni becomes the active/selected tab (was Zaglavlje before)
$wsR.Activate() | Out-Null
